$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.202.56"
$ws.Range("E2").Value = "  -3.79%  "

# Row 3
$ws.Range("D3").Value = "2.264.15"
$ws.Range("E3").Value = "  -4.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.71"
$ws.Range("E5").Value = "  -2.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.37"
$ws.Range("E6").Value = "  -2.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("E8").Value = "  -4.14%  "

# Row 9
$ws.Range("D9").Value = "2.270.68"
$ws.Range("E9").Value = "  -4.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0926"
$ws.Range("E10").Value = "  -5.85%  "

# Row 11
$ws.Range("E11").Value = "  -1.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.74"
$ws.Range("E12").Value = "  +1.91%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.317"
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
$ws.Range("D14").Value = "2.663.50"
$ws.Range("E14").Value = "  -4.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.32"
$ws.Range("E15").Value = "  -1.21%  "

# Row 16
$ws.Range("D16").Value = "54.144.32"
$ws.Range("E16").Value = "  -3.81%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000129"
$ws.Range("E17").Value = "  -2.59%  "

# Row 18
$ws.Range("D18").Value = "2.245.60"
$ws.Range("E18").Value = "  -5.64%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.65"
$ws.Range("E19").Value = "  -4.10%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.98"
$ws.Range("E20").Value = "  -0.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "301.49"
$ws.Range("E21").Value = "  -1.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  -2.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.80"
$ws.Range("E24").Value = "  -1.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.367"
$ws.Range("E26").Value = "  -0.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.144"
$ws.Range("E27").Value = "  -2.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  -2.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.08"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0695"
$ws.Range("E30").Value = "  -2.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.61"
$ws.Range("E31").Value = "  -1.72%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.78"
$ws.Range("E33").Value = "  +0.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("E35").Value = "  -1.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.48"
$ws.Range("E36").Value = "  -0.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.17"
$ws.Range("E37").Value = "  -0.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.850"
$ws.Range("E38").Value = "  +6.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("E39").Value = "  -4.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.73"
$ws.Range("E40").Value = "  -0.78%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -2.09%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.368"
$ws.Range("E42").Value = "  +0.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.32"
$ws.Range("E43").Value = "  -0.68%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.71"
$ws.Range("E44").Value = "  -1.21%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "122.51"
$ws.Range("E45").Value = "  -6.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0882"
$ws.Range("E46").Value = "  -2.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.541"
$ws.Range("E47").Value = "  -4.61%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "238.43"
$ws.Range("E48").Value = "  -1.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0474"
$ws.Range("E49").Value = "  -1.89%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0204"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.28"
$ws.Range("E51").Value = "  -3.50%  "
